{"js": "// Update the title date, then update every math-fact cell in the single\n// 20x5 table with the new values, preserving cell order (row-major, same\n// order as the document XML).\n\nconst newDate = \"2024-05-03 Friday\";\nconst newGrid = [\n  [\"97-94=3\", \"46-13=33\", \"50-33=17\", \"25+47=72\", \"21-9=12\"],\n  [\"66+21=87\", \"99-3=96\", \"42-28=14\", \"68+31=99\", \"6+1=7\"],\n  [\"55-37=18\", \"78-69=9\", \"76-40=36\", \"72-59=13\", \"8+47=55\"],\n  [\"54+33=87\", \"8+63=71\", \"88-67=21\", \"94-29=65\", \"74-10=64\"],\n  [\"49+8=57\", \"5+62=67\", \"98-90=8\", \"49+26=75\", \"30+6=36\"],\n  [\"58-55=3\", \"57+15=72\", \"36+29=65\", \"71-1=70\", \"24+61=85\"],\n  [\"4+71=75\", \"7+13=20\", \"18+16=34\", \"74-63=11\", \"96-73=23\"],\n  [\"11+15=26\", \"26+72=98\", \"77-55=22\", \"67-2=65\", \"41+40=81\"],\n  [\"12+65=77\", \"73-35=38\", \"23+5=28\", \"41-33=8\", \"20+49=69\"],\n  [\"5+55=60\", \"0+48=48\", \"79-10=69\", \"96-52=44\", \"37-19=18\"],\n  [\"0+93=93\", \"85-9=76\", \"88-23=65\", \"40-38=2\", \"45-23=22\"],\n  [\"58-42=16\", \"38+30=68\", \"11+83=94\", \"63-62=1\", \"65+3=68\"],\n  [\"35+17=52\", \"86-63=23\", \"27-13=14\", \"40+14=54\", \"18+12=30\"],\n  [\"90-27=63\", \"57+36=93\", \"70-26=44\", \"87+2=89\", \"57-46=11\"],\n  [\"90-2=88\", \"79-2=77\", \"56+4=60\", \"52-33=19\", \"78-44=34\"],\n  [\"39-2=37\", \"27+56=83\", \"58-6=52\", \"39+29=68\", \"19+40=59\"],\n  [\"68-33=35\", \"79-67=12\", \"40+55=95\", \"66-60=6\", \"39+51=90\"],\n  [\"54-25=29\", \"66-10=56\", \"55-28=27\", \"12+37=49\", \"94-44=50\"],\n  [\"37+53=90\", \"66-34=32\", \"76-66=10\", \"36+20=56\", \"52-7=45\"],\n  [\"96-40=56\", \"49-38=11\", \"77-3=74\", \"6-5=1\", \"37-36=1\"],\n];\n\n// 1) Update the title paragraph (first paragraph in the body) that holds\n//    the date string.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n// Replace the run's text in-place regardless of its current value.\ntitlePara.insertText(newDate, Word.InsertLocation.replace);\n\n// 2) Update every cell of the single table with the new math facts,\n//    preserving row/column order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newGrid;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (first paragraph) holding the date.\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-03 Friday\"\n\n# 2) Update every cell of the single 20x5 table with the new math facts,\n#    preserving row/column order (row-major, matching document order).\n$newGrid = @(\n    ,@(\"97-94=3\", \"46-13=33\", \"50-33=17\", \"25+47=72\", \"21-9=12\")\n    ,@(\"66+21=87\", \"99-3=96\", \"42-28=14\", \"68+31=99\", \"6+1=7\")\n    ,@(\"55-37=18\", \"78-69=9\", \"76-40=36\", \"72-59=13\", \"8+47=55\")\n    ,@(\"54+33=87\", \"8+63=71\", \"88-67=21\", \"94-29=65\", \"74-10=64\")\n    ,@(\"49+8=57\", \"5+62=67\", \"98-90=8\", \"49+26=75\", \"30+6=36\")\n    ,@(\"58-55=3\", \"57+15=72\", \"36+29=65\", \"71-1=70\", \"24+61=85\")\n    ,@(\"4+71=75\", \"7+13=20\", \"18+16=34\", \"74-63=11\", \"96-73=23\")\n    ,@(\"11+15=26\", \"26+72=98\", \"77-55=22\", \"67-2=65\", \"41+40=81\")\n    ,@(\"12+65=77\", \"73-35=38\", \"23+5=28\", \"41-33=8\", \"20+49=69\")\n    ,@(\"5+55=60\", \"0+48=48\", \"79-10=69\", \"96-52=44\", \"37-19=18\")\n    ,@(\"0+93=93\", \"85-9=76\", \"88-23=65\", \"40-38=2\", \"45-23=22\")\n    ,@(\"58-42=16\", \"38+30=68\", \"11+83=94\", \"63-62=1\", \"65+3=68\")\n    ,@(\"35+17=52\", \"86-63=23\", \"27-13=14\", \"40+14=54\", \"18+12=30\")\n    ,@(\"90-27=63\", \"57+36=93\", \"70-26=44\", \"87+2=89\", \"57-46=11\")\n    ,@(\"90-2=88\", \"79-2=77\", \"56+4=60\", \"52-33=19\", \"78-44=34\")\n    ,@(\"39-2=37\", \"27+56=83\", \"58-6=52\", \"39+29=68\", \"19+40=59\")\n    ,@(\"68-33=35\", \"79-67=12\", \"40+55=95\", \"66-60=6\", \"39+51=90\")\n    ,@(\"54-25=29\", \"66-10=56\", \"55-28=27\", \"12+37=49\", \"94-44=50\")\n    ,@(\"37+53=90\", \"66-34=32\", \"76-66=10\", \"36+20=56\", \"52-7=45\")\n    ,@(\"96-40=56\", \"49-38=11\", \"77-3=74\", \"6-5=1\", \"37-36=1\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newGrid.Count; $r++) {\n    $row = $newGrid[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
